# "Change reference style to numbered, not footnote"
#
# With references/citations rendered as numbered (collected) entries
# instead of per-page footnotes, the extra bottom-margin buffer that used
# to reserve room for footnote text at the foot of each page is no longer
# needed, so it is trimmed back down to match the (already small) top
# margin. The caption style is also tightened so a caption's lines always
# stay together on one page/column instead of splitting across a break.

$d = $word.ActiveDocument

# --- Page margins: reclaim the footnote-reservation space at the bottom
#     of the page (w:pgMar w:bottom 1296 -> 288 twips == 64.8pt -> 14.4pt)
$section = $d.Sections(1)
$section.PageSetup.BottomMargin = 14.4

# --- Caption style: keep each caption's lines together (adds <w:keepLines/>
#     to the Caption style's paragraph properties)
$captionStyle = $d.Styles("Caption")
$captionStyle.ParagraphFormat.KeepTogether = $true
